# Update the "dSF" (column F) values for a set of rows with repulled data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    9  = -1
    13 = -1
    14 = -1
    18 = 1
    19 = -3
    21 = -2
    22 = -1
    28 = -5
    30 = -4
    31 = -6
    37 = -1
    40 = -1
    44 = -3
    45 = -2
    46 = -2
    47 = -5
    52 = 1
    53 = -1
    57 = 1
    58 = 2
    61 = 0
    67 = -2
    69 = -4
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
